$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12 (date 44994) - D12 gets new text, wrap text on top of its existing
# font/border formatting.
$ws.Range("D12").WrapText = $true
$ws.Range("D12").Value = "First Day of Joining.`n- Papers Reading`n- Working of basic idea of our design`n- Searching for available equipments"

# Row 13 (date 44995) - D13 gets new text, wrap text added.
$ws.Range("D13").WrapText = $true
$ws.Range("D13").Value = "First Conlusion Compiling`n- Basic Design Proposed`n- Formulas and calculations of Wind, Blades, Shaft, Gears, Dynamo and Battery`n- Major Steps for the fabrication"

# Row 16 (date 44998) - D16 gets new text, wrap text added.
$ws.Range("D16").WrapText = $true
$ws.Range("D16").Value = "Discussion with related persons in automitive department and electrical machines lab regarding dynamos availiblity, performance and market.`nfigured out the the cons in using the car alternators`nWorking on alternate solutions.`nA visit to nearest auto shop for some information regarding dynamos."

# Row 17 (date 44999) - C17 and D17 both get new text, wrap text added.
$ws.Range("C17").WrapText = $true
$ws.Range("C17").Value = "created a methadology report and explained in it the pins utilization.`nAlso created a block diagram of our IOT device system and explained in the methodology report.`nPower system of our project is also explained while creating a block diagram."

$ws.Range("D17").WrapText = $true
$ws.Range("D17").Value = "Working on adding Belt and Pulley system to the design`nAnalysing Data sheets of different manufacturers of alternators.`nListing out the torque and speed requirments of all possibile generators`nCalculating parameters to proof the fesibility of Hub Motor.`nCompile a document of conlusion."

# Update the sheet view to match the new selection/scroll position.
$ws.Range("D19").Select()
$ws.Application.ActiveWindow.ScrollRow = 13
